$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) from the row above (row 10) onto the new row 11
$ws.Range("A10:AQ10").Copy()
$ws.Range("A11:AQ11").PasteSpecial(-4122)

# Fill in the data for the new row (2021 nian)
$ws.Range("A11").Value = "2021年"
$ws.Range("B11").Value = 5259.59
$ws.Range("C11").Value = 1395.47
$ws.Range("D11").Value = 269.4
# E11 intentionally left blank (source cell is an empty string marker)
$ws.Range("F11").Value = 1130.7
$ws.Range("G11").Value = 4054.2
$ws.Range("H11").Value = 408.06
$ws.Range("I11").Value = 3672.06
$ws.Range("J11").Value = 473.03
$ws.Range("K11").Value = 111670.89
$ws.Range("L11").Value = 609.85
$ws.Range("M11").Value = 160.82
$ws.Range("N11").Value = 300.96
$ws.Range("O11").Value = 632.97
$ws.Range("P11").Value = 2620.96
$ws.Range("Q11").Value = 125.34
$ws.Range("R11").Value = 203.42
$ws.Range("S11").Value = 1932.64
$ws.Range("T11").Value = 453.59
$ws.Range("U11").Value = 10852.8
$ws.Range("V11").Value = 281.47
$ws.Range("W11").Value = 3310.94
$ws.Range("X11").Value = 301.96
$ws.Range("Y11").Value = 4294.03
$ws.Range("Z11").Value = 12861.91
$ws.Range("AA11").Value = 532.5700000000001
$ws.Range("AB11").Value = 1597.59
$ws.Range("AC11").Value = 508.48
$ws.Range("AD11").Value = 1186.62
$ws.Range("AE11").Value = 1003.65
$ws.Range("AF11").Value = 29126.22
$ws.Range("AG11").Value = 6600.14
$ws.Range("AH11").Value = 982.34
$ws.Range("AI11").Value = 506.48
$ws.Range("AJ11").Value = 291.35
$ws.Range("AK11").Value = 3156.25
$ws.Range("AL11").Value = 3242.74
$ws.Range("AM11").Value = 3237.62
$ws.Range("AN11").Value = 97.48999999999999
$ws.Range("AO11").Value = 1276.71
$ws.Range("AP11").Value = 2340.59
$ws.Range("AQ11").Value = 377.42
